$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (Location, Site, Exposure period, Notes, Exist)
$data = @(
    @("Abbotsford", "Bodriggy Brewing Company  245 Johnston Street, Abbotsford VIC 3067", "28/12/20 2:50pm-5:30pm", "Case dined at venue", "new"),
    @("Black Rock", "Smile Buffalo Thai restaurant  305 Beach Road, Black Rock VIC 3193", "21/12/20 6:30pm-9:40pm", "Case dined in restaurant", "old"),
    @("Black Rock", "Smile Buffalo Thai restaurant  305 Beach Road, Black Rock VIC 3193", "21/12/20 7:30pm-9:30pm", "Case dined in restaurant", "old"),
    @("Hampton", "Merrymen Cafe, 2 Small Street, Hampton VIC", "28-12-2020 1:30pm-2:30pm", "Case ate in store", "new"),
    @("Keysborough", "Sikh Temple Keysborough, 198-206 Perry Road, Keysborough", "1/01/21 3:00pm-5:00pm", "Case visited venue", "new"),
    @("Springvale", "IKEA Springvale - Cafe and Restaurant, 917 Princes Hwy", "30/12/20 5:30pm - 6.30pm", "Case dined at cafe", "new"),
    @("Springvale", "IKEA Springvale, 917 Princes Hwy", "29/12/20 5:30pm - 6.30pm", "Case shopped at store and dined at cafe", "old")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $row++
}

# Autofit columns to content, matching the bestFit behaviour of the original
$ws.Range("A:E").AutoFit() | Out-Null

# Select entire columns A:E as the final selection state
$ws.Range("A1:E1048576").Select() | Out-Null
